# Updated fitting parameters. Ready to run detail tests.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Update the fitted parameter values
$ws.Range("J2").Value = 0.01525
$ws.Range("K2").Value = 0.1115

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("K3").Select()
